$d = $word.ActiveDocument

# The document currently ends with:
#   ...I implemented fully the sample codes from the tutorial website
#   <empty paragraph carrying a leftover sz=40 run-properties mark>
# We replace that trailing empty paragraph with the new log entries,
# using raw WordprocessingML so every run / proofErr / hyperlink shell
# matches exactly, then wire up the real hyperlink relationship via
# Hyperlinks.Add (so Word mints the relationship + Hyperlink style run
# the normal COM way instead of us faking an r:id).

$w = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:r>
    <w:t xml:space="preserve">Now I\u2019m going to figure out how to generate movie recommendations based on more variables. I reckon I could use the </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>tfidf</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t xml:space="preserve"> to transform all the text attributes and in combination with the normalized numerical values, I could average </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:u w:val="single"/>
    </w:rPr>
    <w:t>ALL</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> of them and then use cosine similarity to generate recommendation.</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:r>
    <w:lastRenderedPageBreak/>
    <w:t xml:space="preserve">I found a great site for learning about TF-IDF working. ( </w:t>
  </w:r>
  <w:r>
    <w:t>https://www.geeksforgeeks.org/machine-learning/understanding-tf-idf-term-frequency-inverse-document-frequency/</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> )</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
'@

$xml = $xml.Replace("\u2019", [char]0x2019)

$target = $d.Paragraphs.Last.Range
$target.InsertXML($xml)

# Turn the plain-text URL we just inserted into a real hyperlink,
# the same way Word does it interactively: select the display text,
# then Hyperlinks.Add wraps it in a w:hyperlink + mints the relationship.
$linkRange = $d.Content
$linkRange.Find.ClearFormatting()
$linkRange.Find.Execute("https://www.geeksforgeeks.org/machine-learning/understanding-tf-idf-term-frequency-inverse-document-frequency/", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
$d.Hyperlinks.Add($linkRange, "https://www.geeksforgeeks.org/machine-learning/understanding-tf-idf-term-frequency-inverse-document-frequency/", "", "", $linkRange.Text) | Out-Null

Write-Output "done"
